{"js": "const body = context.document.body;\n\n// Each entry is an exact, uniquely-occurring substring in the document and\n// its replacement text. Using body.search + Range.insertText(\"Replace\")\n// lets us target the precise text regardless of how it is split across\n// <w:r> runs in the underlying OOXML.\nconst edits = [\n  {\n    find: \"during the 1990s, Netscape Corporation\",\n    replace: \"during the late 1990s, Netscape Corporation\"\n  },\n  {\n    find: \"Microsoft, Netscape, Sun Microsystems all contributed\",\n    replace: \"Microsoft, Netscape, and Sun Microsystems all contributed\"\n  },\n  {\n    find: \"the majority status it does today, and was rightfully cautious of implementing agreed upon standard first\",\n    replace: \"the majority status that it does today, and Microsoft was rightfully cautious of implementing agreed upon standards first\"\n  },\n  {\n    find: \"Were this to happen, Internet Explorer may not be able to view many web pages\",\n    replace: \"Were this to happen, Internet Explorer would not be able to be used to view many web pages\"\n  },\n  {\n    find: \"This has proved overall to be beneficial for the web; even to this day\",\n    replace: \"This has proved overall to be beneficial for the Web; even to this day\"\n  },\n  {\n    find: \"This would have indeed been a tragedy because Web developers would be much more difficult\",\n    replace: \"This would have indeed been a tragedy because Web development would be much more difficult\"\n  }\n];\n\nfor (const { find, replace } of edits) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"expected exactly 1 match for \" + JSON.stringify(find) +\n      \" but found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# List of (find, replace) pairs to apply across the document. Each find\n# string is an exact, uniquely-occurring substring, so a single\n# wdReplaceOne pass against $d.Content is sufficient and safe for each.\n$edits = @(\n  @{ Find = \"during the 1990s, Netscape Corporation\"; Replace = \"during the late 1990s, Netscape Corporation\" },\n  @{ Find = \"Microsoft, Netscape, Sun Microsystems all contributed\"; Replace = \"Microsoft, Netscape, and Sun Microsystems all contributed\" },\n  @{ Find = \"the majority status it does today, and was rightfully cautious of implementing agreed upon standard first\"; Replace = \"the majority status that it does today, and Microsoft was rightfully cautious of implementing agreed upon standards first\" },\n  @{ Find = \"Were this to happen, Internet Explorer may not be able to view many web pages\"; Replace = \"Were this to happen, Internet Explorer would not be able to be used to view many web pages\" },\n  @{ Find = \"This has proved overall to be beneficial for the web; even to this day\"; Replace = \"This has proved overall to be beneficial for the Web; even to this day\" },\n  @{ Find = \"This would have indeed been a tragedy because Web developers would be much more difficult\"; Replace = \"This would have indeed been a tragedy because Web development would be much more difficult\" }\n)\n\nforeach ($edit in $edits) {\n  $rng = $d.Content\n  $find = $rng.Find\n  $find.Text = $edit.Find\n  $find.Replacement.Text = $edit.Replace\n  $find.Forward = $true\n  $find.Wrap = 0\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.MatchSoundsLike = $false\n  $find.MatchAllWordForms = $false\n\n  $ok = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 0, $false, $edit.Replace, 2)\n  if (-not $ok) {\n    throw \"Find/Replace failed for: $($edit.Find)\"\n  }\n}\n\nWrite-Output \"done\"\n"}
